# Daily attendance processing - 2026-01-11 15:34:20
# Normalizes the "Recorded By" (column G) value ordering on the
# "Session Analysis Results" sheet:
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System, system" -> "backup@backdoor.com, system, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value2 = "backup@backdoor.com, system, System"
    }
}
